$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.180.18"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.237.22"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'293.60"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'88.80"
$ws.Range("E6").Value = "  +6.06%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'31.21"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "2.580.46"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "'14.23"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "2.278.90"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "'0.739"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "40.131.97"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "'11.68"
$ws.Range("E20").Value = "  +11.67%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'5.88"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'66.25"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'2.49"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").Value = "'1.86"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'9.36"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'33.13"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").Value = "'152.67"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "'5.01"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("D38").Value = "'16.32"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'0.113"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").Value = "2.117.63"
$ws.Range("E42").Value = "  +9.90%  "
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("E44").Value = "  +6.65%  "
$ws.Range("D45").Value = "'18.44"
$ws.Range("E45").Value = "  +11.05%  "
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("D47").Value = "'10.12"
$ws.Range("E47").Value = "  +10.56%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "2.450.31"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'71.47"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'1.48"
$ws.Range("E51").Value = "  +6.49%  "

# Force numeric-looking price cells to remain plain text (no explicit "Text" format),
# matching the original inlineStr cells which carry no style override.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
